$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Old row 2 ("Hiver"/"Eté"/"Année" sub-header) is removed entirely; rows below shift up by one.
$ws.Rows.Item(2).Delete()

# Rewrite row 1 as the new header row (idx, idx2, Name, Date Start, Date End, (m3/s), (MW1), (MW2), (GWh) Winter/Summer/Year)
$ws.Range("A1:K1").ClearContents()
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# F1:K1 use a distinct (font-only) cell style -- same font as the data columns (Arial 9)
# but without a number-format override. Go through a transient named style so the
# engine emits a plain cellXf (applyFont only) instead of reusing/augmenting existing ones,
# then delete the named style so only the direct formatting (cellXfs entry) remains.
$st = $wb.Styles.Add("HeaderStyleTmp")
$st.Font.Name = "Arial"
$st.Font.Size = 9
$ws.Range("F1:K1").Style = "HeaderStyleTmp"
$wb.Styles.Item("HeaderStyleTmp").Delete()
